# Inserts a new weekly price record at row 147 of the "Espinaca" sheet,
# pushing the existing rows 147-174 down to 148-175.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 147:174 down to 148:175, inserting a blank row at 147
# (copies formatting from the row above, same as Excel's default Insert).
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new record.
$ws.Range("A147").Value = 8
$ws.Range("B147").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C147").Value = 'Coquimbo'
$ws.Range("D147").Value = 44505
$ws.Range("E147").Value = 4
$ws.Range("F147").Value = 100112012
$ws.Range("G147").Value = 'Espinaca'
$ws.Range("H147").Value = 'Sin especificar'
$ws.Range("I147").Value = 'Primera'
$ws.Range("J147").Value = 3200
$ws.Range("K147").Value = 400
$ws.Range("L147").Value = 500
$ws.Range("M147").Value = 450
$ws.Range("N147").Value = '$/atado 300 a 500 gramos'
$ws.Range("O147").Value = 'Provincia del Elquí'
$ws.Range("P147").Value = 900
$ws.Range("Q147").Value = 0.5
$ws.Range("R147").Value = 'Hortaliza'

# Give D147 the same date style ("s=2", i.e. the date number format) as the
# rest of the date column.
$ws.Range("D147").NumberFormat = $ws.Range("D148").NumberFormat
